$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 438
$ws.Range("I98").Value = 438
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 438
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 1060
$ws.Range("N98").Value = $null
$ws.Range("H100").Value = 1398.4286
$ws.Range("I100").Value = 1298.1666
$ws.Range("K100").Value = 1298.1666
$ws.Range("M100").Value = -757.1666
$ws.Range("H106").Value = 3212
$ws.Range("I106").Value = 3212
$ws.Range("K106").Value = 3212
$ws.Range("M106").Value = -2581
$ws.Range("H122").Value = 438
$ws.Range("I122").Value = 438
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1314
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 1136
$ws.Range("N122").Value = $null
$ws.Range("H132").Value = 17796.916
$ws.Range("I132").Value = 17796.916
$ws.Range("K132").Value = 53390.74800000001
$ws.Range("M132").Value = -50860.74800000001
$ws.Range("H135").Value = 605.7857
$ws.Range("I135").Value = 605.7857
$ws.Range("K135").Value = 5452.071300000001
$ws.Range("M135").Value = -2917.071300000001
$ws.Range("H137").Value = 2555.327
$ws.Range("I137").Value = 1622.4642
$ws.Range("J137").Value = 3643.6667
$ws.Range("K137").Value = 4867.392599999999
$ws.Range("L137").Value = 10931.0001
$ws.Range("M137").Value = -2317.392599999999
$ws.Range("N137").Value = -16031.0001
$ws.Range("H138").Value = 3257.6
$ws.Range("J138").Value = 3175.8
$ws.Range("L138").Value = 9527.400000000001
$ws.Range("N138").Value = -19807.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2514.3076
$ws.Range("I45").Value = 2163.625
$ws.Range("J45").Value = 3075.4
$ws.Range("K45").Value = 2163.625
$ws.Range("L45").Value = 3075.4
$ws.Range("M45").Value = -1786.625
$ws.Range("N45").Value = -3829.4
$ws.Range("H56").Value = 14517.5
$ws.Range("J56").Value = 14656.667
$ws.Range("L56").Value = 14656.667
$ws.Range("N56").Value = -16140.667
$ws.Range("H61").Value = 3230.2
$ws.Range("I61").Value = 3230.2
$ws.Range("K61").Value = 3230.2
$ws.Range("M61").Value = -3018.2
$ws.Range("H97").Value = 1278.8125
$ws.Range("I97").Value = 875
$ws.Range("K97").Value = 875
$ws.Range("M97").Value = -379
$ws.Range("H102").Value = 3995.611
$ws.Range("I102").Value = 2794.7334
$ws.Range("K102").Value = 2794.7334
$ws.Range("M102").Value = -1172.7334
$ws.Range("H122").Value = 4199.8335
$ws.Range("I122").Value = 3999.6667
$ws.Range("J122").Value = 4400
$ws.Range("K122").Value = 11999.0001
$ws.Range("L122").Value = 13200
$ws.Range("M122").Value = -9549.000100000001
$ws.Range("N122").Value = -18100
$ws.Range("H136").Value = 3230.2
$ws.Range("I136").Value = 3230.2
$ws.Range("K136").Value = 9690.599999999999
$ws.Range("M136").Value = -7140.599999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1055.3334
$ws.Range("I80").Value = 771.6
$ws.Range("J80").Value = 1410
$ws.Range("K80").Value = 771.6
$ws.Range("L80").Value = 1410
$ws.Range("M80").Value = 226.4
$ws.Range("N80").Value = -3406
$ws.Range("H83").Value = 1055.3334
$ws.Range("I83").Value = 771.6
$ws.Range("J83").Value = 1410
$ws.Range("K83").Value = 3858
$ws.Range("L83").Value = 7050
$ws.Range("M83").Value = 1134
$ws.Range("N83").Value = -17034
$ws.Range("H86").Value = 8450
$ws.Range("I86").Value = 6750
$ws.Range("J86").Value = 9583.333000000001
$ws.Range("K86").Value = 6750
$ws.Range("L86").Value = 9583.333000000001
$ws.Range("M86").Value = -5627
$ws.Range("N86").Value = -11829.333
$ws.Range("H89").Value = 8450
$ws.Range("I89").Value = 6750
$ws.Range("J89").Value = 9583.333000000001
$ws.Range("K89").Value = 33750
$ws.Range("L89").Value = 47916.665
$ws.Range("M89").Value = -28134
$ws.Range("N89").Value = -59148.665
$ws.Range("H94").Value = 1809.909
$ws.Range("J94").Value = 1603
$ws.Range("L94").Value = 1603
$ws.Range("N94").Value = -2505
$ws.Range("H99").Value = 2175.6428
$ws.Range("I99").Value = 2225.7144
$ws.Range("K99").Value = 2225.7144
$ws.Range("M99").Value = -727.7143999999998
$ws.Range("H134").Value = 5000
$ws.Range("I134").Value = 5000
$ws.Range("K134").Value = 15000
$ws.Range("M134").Value = -12465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 5000
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 5000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 5000
$ws.Range("M6").Value = $null
$ws.Range("N6").Value = -5226

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 162.5
$ws.Range("I38").Value = 122
$ws.Range("K38").Value = 366
$ws.Range("M38").Value = -19
$ws.Range("H139").Value = 2406.0557
$ws.Range("I139").Value = 1613.9333
$ws.Range("K139").Value = 4841.7999
$ws.Range("M139").Value = 298.2001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3172
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 3688
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 3688
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -5684
$ws.Range("H83").Value = 3172
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 3688
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 18440
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -28424

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2573.5715
$ws.Range("I22").Value = 2333.3333
$ws.Range("J22").Value = 2753.75
$ws.Range("K22").Value = 2333.3333
$ws.Range("L22").Value = 2753.75
$ws.Range("M22").Value = -2038.3333
$ws.Range("N22").Value = -3343.75
$ws.Range("H27").Value = 2573.5715
$ws.Range("I27").Value = 2333.3333
$ws.Range("J27").Value = 2753.75
$ws.Range("K27").Value = 2333.3333
$ws.Range("L27").Value = 2753.75
$ws.Range("M27").Value = -2226.3333
$ws.Range("N27").Value = -2967.75
$ws.Range("H40").Value = 5375
$ws.Range("I40").Value = 5375
$ws.Range("K40").Value = 5375
$ws.Range("M40").Value = -5239
$ws.Range("H93").Value = 896.6667
$ws.Range("I93").Value = 845
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 845
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = 403
$ws.Range("N93").Value = -3496

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 20000000
$ws.Range("J5").Value = 20000000
$ws.Range("L5").Value = 20000000
$ws.Range("N5").Value = -20000224
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").Value = $null
$ws.Range("H126").Value = 6213.7856
$ws.Range("I126").Value = 3598.6
$ws.Range("K126").Value = 10795.8
$ws.Range("M126").Value = -8325.799999999999
$ws.Range("H132").Value = 2162.5386
$ws.Range("I132").Value = 2147.1904
$ws.Range("J132").Value = 2227
$ws.Range("K132").Value = 6441.5712
$ws.Range("L132").Value = 6681
$ws.Range("M132").Value = -3911.5712
$ws.Range("N132").Value = -11741
